# Apply crypto price/volume updates from the Tue May 14 04:34:48 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume 1h) columns hold text-formatted numbers/percentages
# (leading/trailing spaces, thousands dots, scientific-notation-prone decimals).
# Force the Text number format before writing so Excel's COM value-coercion
# doesn't reinterpret them as actual numbers (which would mangle the string).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.572.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.940.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.938.28"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.58%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.35"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.429.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.591.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.941.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.663"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.18"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.94"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.14%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.58"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.51%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +20.42%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.10"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.988"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.64"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.36"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.43"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.696.88"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0339"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "352.82"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.40%  "
